# update: bereinigtes Abschlussquote-Excel & aktualisiertes Profil-Notebook
#
# Cleans up the sheet names (fixes a stray trailing dot and normalizes
# "data"/"dict" -> "Data"/"Dict" casing) and restores each sheet's last
# active-cell selection.

$wb = $excel.ActiveWorkbook

# --- Fix / normalize worksheet names -------------------------------------
$wb.Worksheets.Item(4).Name = "T2_SekII_1st_25_Kant_Dict"
$wb.Worksheets.Item(5).Name = "T3_Matura_Merkm_Data"
$wb.Worksheets.Item(6).Name = "T3_Matura_Merk_Dict"
$wb.Worksheets.Item(7).Name = "T4_Matura_Kant_Data"
$wb.Worksheets.Item(8).Name = "T4_Matura_Kant_Dict"

# --- Restore per-sheet selections -----------------------------------------
# Sheets not listed here (1 and 3) keep their existing selection untouched.

$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("C15:C16").Select()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$ws5.Range("G23").Select()

$ws6 = $wb.Worksheets.Item(6)
$ws6.Activate()
$ws6.Range("B25").Select()

$ws7 = $wb.Worksheets.Item(7)
$ws7.Activate()
$ws7.Range("I15").Select()

$ws8 = $wb.Worksheets.Item(8)
$ws8.Activate()
$ws8.Range("J20").Select()

# Sheet 2 (T1_SekII_1st_25_Merkm_Dict) is the sheet that stays active/selected
# in the workbook (tabSelected / activeTab) once the file is saved again.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("C19").Select()
